$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.129.60'
$ws.Range('E2').Value = '  +0.86%  '
$ws.Range('D3').Value = '2.267.67'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '305.69'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '95.87'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.12%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.490'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.97%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '35.23'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +8.48%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0793'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.48%  '
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.64'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.46%  '
$ws.Range('D14').Value = '2.620.87'
$ws.Range('E15').Value = '  +0.99%  '
$ws.Range('D16').Value = '2.276.12'
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.794'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.44%  '
$ws.Range('D18').Value = '42.046.67'
$ws.Range('E18').Value = '  +0.82%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.45'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.71%  '
$ws.Range('D20').Value = '0.0₃0905'
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.98'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.22%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '67.84'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.34%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '237.83'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.70%  '
$ws.Range('E24').Value = '  -0.46%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.93'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.42%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '23.71'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.28%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '36.88'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +5.49%  '
$ws.Range('E29').Value = '  -0.77%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.10'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.59%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '159.97'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.24'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.33%  '
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.18'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +4.88%  '
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '17.06'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.08%  '
$ws.Range('E37').Value = '  +0.29%  '
$ws.Range('E38').Value = '  -0.29%  '
$ws.Range('E39').Value = '  +1.98%  '
$ws.Range('E40').Value = '  -1.40%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.03'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.66%  '
$ws.Range('E42').Value = '  +8.22%  '
$ws.Range('D43').Value = '1.987.27'
$ws.Range('E43').Value = '  -1.11%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '19.06'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -3.20%  '
$ws.Range('E45').Value = '  +0.50%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '9.92'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.98%  '
$ws.Range('E47').Value = '  +0.56%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '53.11'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.48%  '
$ws.Range('E49').Value = '  +0.49%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '72.01'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '90.98'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.80%  '
